$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")
$ws.Activate()

$range = $ws.Range("F1:F151")
$range.Select()
$range.ClearContents()
